$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17, pushing existing rows 17-19 down to 18-20
$ws.Rows.Item(17).Insert()

# New row 17 values: same as old row 17, except D (date) and J (volumen) change
$ws.Cells.Item(17, 1).Value = 8
$ws.Cells.Item(17, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 45007
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 100112039
$ws.Cells.Item(17, 7).Value = "Ciboulette"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 1160
$ws.Cells.Item(17, 11).Value = 2000
$ws.Cells.Item(17, 12).Value = 2500
$ws.Cells.Item(17, 13).Value = 2250
$ws.Cells.Item(17, 14).Value = "`$/docena de atados"
$ws.Cells.Item(17, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(17, 16).Value = 750
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = "Hortaliza"
